$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap values in columns I, K, M, O for data rows 2..25 ---
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2 (was 1)
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1 (was 2)
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2 (was 1)
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1 (was 2)
}

# --- Add new columns P (16) and Q (17) ---
# Header row (row 1): continue the sequential numbering 14, 15,
# matching the formatting (bold, thin border on all sides,
# centered horizontally, top vertical alignment) used by the rest of row 1.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# Data rows 2..25 for columns P and Q, value 2
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 2
    $ws.Cells.Item($r, 17).Value = 2
}
